# Weekly update: insert a new price record for "Vega Modelo de Temuco" -
# Jengibre (row 151), pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 151; this shifts rows
# 151..184 down to 152..185 and keeps their values/formatting intact.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with this week's record.
$ws.Range("A151").Value = 10
$ws.Range("B151").Value = "Vega Modelo de Temuco"
$ws.Range("C151").Value = "La Araucanía"
$ws.Range("D151").Value = 44785
$ws.Range("E151").Value = 9
$ws.Range("F151").Value = 100114007
$ws.Range("G151").Value = "Jengibre"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 20
$ws.Range("K151").Value = 16000
$ws.Range("L151").Value = 16000
$ws.Range("M151").Value = 16000
$ws.Range("N151").Value = "$/caja 13 kilos"
$ws.Range("O151").Value = "Perú"
$ws.Range("P151").Value = 1231
$ws.Range("Q151").Value = 13
$ws.Range("R151").Value = "Hortaliza"
